$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.266.91'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '1.659.20'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.79%  '
$ws.Range("D5").Value = '219.71'
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("D6").Value = '0.5248'
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("D8").Value = '0.2677'
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").Value = '0.06378'
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").Value = '20.73'
$ws.Range("E10").Value = '  -1.43%  '
$ws.Range("D11").Value = '0.07765'
$ws.Range("E11").Value = '  -1.11%  '
$ws.Range("D12").Value = '4.586'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '1.575.74'
$ws.Range("E13").Value = '  -5.84%  '
$ws.Range("D14").Value = '1.886.03'
$ws.Range("E14").Value = '  -0.93%  '
$ws.Range("D15").Value = '0.5691'
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").Value = '0.0₅8202'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '65.61'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '26.250.28'
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").Value = '4.729'
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").Value = '192.40'
$ws.Range("E21").Value = '  -3.40%  '
$ws.Range("D22").Value = '10.40'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").Value = '6.049'
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").Value = '  -0.79%  '
$ws.Range("D25").Value = '143.56'
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("D26").Value = '0.1207'
$ws.Range("E26").Value = '  -2.29%  '
$ws.Range("D27").Value = '7.290'
$ws.Range("E27").Value = '  +0.39%  '
$ws.Range("D28").Value = '16.03'
$ws.Range("E28").Value = '  -1.58%  '
$ws.Range("D29").Value = '1.493'
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").Value = '0.05647'
$ws.Range("E30").Value = '  -4.13%  '
$ws.Range("D31").Value = '1.282'
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").Value = '3.507'
$ws.Range("E32").Value = '  -2.26%  '
$ws.Range("D33").Value = '3.381'
$ws.Range("E33").Value = '  +1.93%  '
$ws.Range("D34").Value = '1.588'
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("D35").Value = '2.807'
$ws.Range("E35").Value = '  -1.62%  '
$ws.Range("D36").Value = '0.9485'
$ws.Range("E36").Value = '  -2.40%  '
$ws.Range("E37").Value = '  -1.42%  '
$ws.Range("D38").Value = '0.5780'
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("D39").Value = '0.01602'
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("D40").Value = '5.923'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").Value = '2.578'
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("D44").Value = '1.031.67'
$ws.Range("E44").Value = '  -4.76%  '
$ws.Range("D45").Value = '102.44'
$ws.Range("E45").Value = '  -1.84%  '
$ws.Range("D46").Value = '1.796.29'
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("D47").Value = '58.66'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").Value = '0.0₈107'
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("D49").Value = '1.005'
$ws.Range("E49").Value = '  -0.94%  '
$ws.Range("D50").Value = '0.05313'
$ws.Range("E50").Value = '  +2.72%  '
$ws.Range("D51").Value = '8.033'
$ws.Range("E51").Value = '  -0.26%  '
